# Update "想去人数" (interest/attendee count) values across sheets,
# matching the regenerated gh-pages data output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$wsExpo  = $wb.Worksheets.Item("展览")
$wsPerf  = $wb.Worksheets.Item("演出")
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsAll   = $wb.Worksheets.Item("全部类型")

# 展览 (Exhibition)
$wsExpo.Range("F6").Value  = 168
$wsExpo.Range("F8").Value  = 315
$wsExpo.Range("F11").Value = 28478
$wsExpo.Range("F12").Value = 3383
$wsExpo.Range("F20").Value = 608
$wsExpo.Range("F25").Value = 46

# 演出 (Performance)
$wsPerf.Range("F16").Value = 44

# 本地生活 (Local Life)
$wsLocal.Range("F3").Value = 251
$wsLocal.Range("F4").Value = 1164

# 全部类型 (All Types) - aggregate sheet mirrors the above events
$wsAll.Range("F3").Value  = 251
$wsAll.Range("F4").Value  = 1164
$wsAll.Range("F12").Value = 168
$wsAll.Range("F14").Value = 315
$wsAll.Range("F25").Value = 3383
$wsAll.Range("F27").Value = 44
$wsAll.Range("F28").Value = 44
$wsAll.Range("F34").Value = 608
$wsAll.Range("F38").Value = 46
